# Applies:
#  - Merge the "Precision Forward " + "Movement" runs (which straddle the
#    _GoBack bookmark) in the "Module 3" paragraph into a single run
#    "Precision Forward Movement", then split a new paragraph off after it
#    (ind left=360, no list numbering) that holds the _GoBack bookmark and a
#    single-space run.
#  - Delete the "Module 4" .. "Module 8" list paragraphs entirely.

$d = $word.ActiveDocument

# --- Step 1: locate "Precision Forward " and "Movement" (they are adjacent,
# separated only by the collapsed _GoBack bookmark) ---
$rForward = $d.Content
$rForward.Find.ClearFormatting()
$foundForward = $rForward.Find.Execute("Precision Forward ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$rMovement = $d.Content
$rMovement.Find.ClearFormatting()
$foundMovement = $rMovement.Find.Execute("Movement", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $foundForward -or -not $foundMovement) {
    throw "Could not locate 'Precision Forward ' / 'Movement' runs"
}

# --- Step 2: merge them into a single run "Precision Forward Movement",
# replacing only that span so the "Module 3 - " prefix (and its exact
# en-dash run) is left completely untouched. ---
$mergeRange = $d.Range($rForward.Start, $rMovement.End)
$mergeXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Precision Forward Movement</w:t></w:r></w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'
$mergeRange.InsertXML($mergeXml)

# --- Step 3: split a new paragraph off right after the merged run (i.e.
# right before the paragraph mark of the "Module 3" paragraph). ---
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -eq "Module 3 – Precision Forward Movement`r" -or
        $cand.Range.Text -like "Module 3*Precision Forward Movement`r") {
        $targetPara = $cand
        break
    }
}
if (-not $targetPara) {
    throw "Could not find the merged 'Module 3' paragraph"
}

$splitPos = $targetPara.Range.End - 1
$splitRange = $d.Range($splitPos, $splitPos)
$splitRange.InsertParagraphAfter()

$newParaIndex = $targetPara.Index + 1
$newPara = $d.Paragraphs.Item($newParaIndex)

# --- Step 4: give the new paragraph its own formatting (ind left=360, no
# list numbering/style) and content (the _GoBack bookmark plus a single
# space run). ---
$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
    '<w:pPr><w:ind w:left="360"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '</w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($newParaXml)

# --- Step 5: delete the "Module 4" .. "Module 8" paragraphs entirely. ---
$firstIdx = $null
$lastIdx = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "Module 4*" -or $txt -like "Module 5*" -or $txt -like "Module 6*" -or
        $txt -like "Module 7*" -or $txt -like "Mod*ule 8*") {
        if ($null -eq $firstIdx) { $firstIdx = $i }
        $lastIdx = $i
    }
}

if ($firstIdx -ne $null) {
    $delStart = $d.Paragraphs.Item($firstIdx).Range.Start
    $delEnd = $d.Paragraphs.Item($lastIdx).Range.End
    $delRange = $d.Range($delStart, $delEnd)
    $delRange.Delete()
}

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Output ("[$i] '" + $d.Paragraphs.Item($i).Range.Text + "'")
}
